# MBA Options workbook update
# - Drop the "University Ranking" column (was C) - its free-text ranking
#   paragraph is replaced by four dedicated ranking columns at the end.
# - Rename "Website URL" (now shifted to column D) to "Program Link".
# - Add four new trailing columns: Webometrics World, Webometrics National,
#   US News, QS Ranking - with per-row numeric/NA values.
# - Move the yellow "needs attention" highlight from the link column to the
#   Program Name column (same set of highlighted rows as before).
# - Shrink the row height on the three data rows now that the long wrapped
#   ranking paragraph is gone.
# - Refresh the _FilterDatabase defined name to the new A1:K43 extent.
# - Nudge the remembered selection the way the author last left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the "University Ranking" column (column C).
# ---------------------------------------------------------------------------
$ws.Columns("C").Delete()

# ---------------------------------------------------------------------------
# 2. Rename the (now shifted) "Website URL" header in column D.
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = "Program Link"

# ---------------------------------------------------------------------------
# 3. Move the yellow highlight from column D (Program Link) to column C
#    (Program Name) for the rows that were previously flagged.
# ---------------------------------------------------------------------------
$highlightRows = @(2,3,4,5,6,8,14,15,16,17,18,23,24,25,27,28,29,30,31,32,36,37,43)
foreach ($r in $highlightRows) {
    $ws.Range("D$r").Interior.Color = 16777215
    $ws.Range("C$r").Interior.Color = 65535
}

# ---------------------------------------------------------------------------
# 4. Populate the new ranking columns (L:O) for the header and the three
#    data rows, writing left-to-right so each new column inherits the
#    correct row-level formatting as it is created.
# ---------------------------------------------------------------------------
$ws.Range("L1").Value = "Webometrics World"
$ws.Range("M1").Value = "Webometrics National"
$ws.Range("N1").Value = "US News"
$ws.Range("O1").Value = "QS Ranking"

$rankingData = @{
    2 = @(23, 5,  "NA", 78)
    3 = @(45, 6,  77,   435)
    4 = @(43, 34, 34,   45)
}
foreach ($r in $rankingData.Keys) {
    $vals = $rankingData[$r]
    $ws.Range("L$r").Value = $vals[0]
    $ws.Range("M$r").Value = $vals[1]
    $ws.Range("N$r").Value = $vals[2]
    $ws.Range("O$r").Value = $vals[3]
}

# Match the look of the existing header / body cells for the new columns.
$headerRange = $ws.Range("L1:O1")
$headerRange.Font.Bold = $true
$headerRange.Font.Color = 16777215
$headerRange.Interior.Color = 6567712
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4131
$headerRange.VerticalAlignment = -4108
$headerRange.WrapText = $true

$bodyRange = $ws.Range("L2:O43")
$bodyRange.Interior.Color = 16777215
$bodyRange.Borders.LineStyle = 1
$bodyRange.HorizontalAlignment = -4131
$bodyRange.VerticalAlignment = -4160
$bodyRange.WrapText = $true

# Match column width of the adjoining "Yearly Tuition Fees" column (K).
$ws.Range("L1:O43").ColumnWidth = $ws.Range("K1").ColumnWidth()

# ---------------------------------------------------------------------------
# 5. The wrapped ranking paragraph is gone, so the three data rows no longer
#    need the tall row height.
# ---------------------------------------------------------------------------
$ws.Rows("2:4").RowHeight = 28.8

# ---------------------------------------------------------------------------
# 6. Refresh the hidden _FilterDatabase name to the new used range.
# ---------------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -like "*FilterDatabase*") {
        $n.RefersTo = "=Sheet!`$A`$1:`$K`$43"
    }
}

# ---------------------------------------------------------------------------
# 7. Leave the selection where the author last left it.
# ---------------------------------------------------------------------------
$ws.Range("D14").Select()
